$wb = $excel.ActiveWorkbook

# --- Sheet "Summary": update A5 value and the selected cell/range ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A5").Value = 1
$wsSummary.Activate()
$wsSummary.Range("F5").Select()

# --- Sheet "Day by Day": update ELN values in column E and the selected cell/range ---
$wsDayByDay = $wb.Worksheets.Item("Day by Day")
$wsDayByDay.Range("E5").Value = "00716727-0068"
$wsDayByDay.Range("E6").Value = "00716727-0069"
$wsDayByDay.Range("E7").Value = "00716727-0069"
$wsDayByDay.Activate()
$wsDayByDay.Range("G11").Select()
